$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.863.33'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.14%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.077.81'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.05%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.46%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.14'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.15%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.076.49'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.91%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.514'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.11%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.40'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.19%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.150'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.72%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.471'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.32%  '

$ws.Range("E13").Value = '  -2.25%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.24'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.34%  '

$ws.Range("E15").Value = '  -2.20%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.584.19'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.15%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.771.22'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.31%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.00'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.39%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.076.15'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.27%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.74%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '487.76'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.25%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.74'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.30%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.686'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.40%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.63'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.39%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.84'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.12%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.23'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.65%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.24'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.94%  '

$ws.Range("E28").Value = '  +0.12%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.75'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.54%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.31'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.26%  '

$ws.Range("E31").Value = '  -1.16%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.81'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.74%  '

$ws.Range("E33").Value = '  -1.69%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0906'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.16%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.19%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.70'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.51%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.954'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.31%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '46.35'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.35%  '

$ws.Range("E39").Value = '  +0.76%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.98'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.82%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.302'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.52%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.32'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.770.95'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.88%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '370.18'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.33%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0345'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.57%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '135.49'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.49%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.46'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.16%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '24.36'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.39%  '

$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.15'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.26%  '

$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.106'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.41%  '
